# Add current price column (現價) as new column B, shifting existing
# columns C..Q to D..R, then replace the two sample rows with index rows
# (TWII / TWOII) and append two new fully-populated stock rows
# (2427.TW and 1264.TWO, the latter an OTC/TWO listing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at B; this shifts 日主力..負加總 from C:Q to D:R
#    and preserves per-cell formatting of those columns automatically.
$ws.Columns("B:B").Insert()

# 2) New header cell for the inserted column - copy the header style
#    (bold + border + centered) from the neighbouring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "現價"

# 3) Row 2 becomes the TWII (加權指數) index row: ticker + price only,
#    the rest of the row's indicator columns are cleared.
$ws.Range("A2").Value = "TWII"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "18,096.07"
$ws.Range("C2:R2").ClearContents()

# 4) Row 3 becomes the TWOII (櫃買指數) index row: ticker + price only.
$ws.Range("A3").Value = "TWOII"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "238.67"
$ws.Range("C3:R3").ClearContents()

# 5) New row 4: 2427.TW, full indicator row.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "2427.TW"
foreach ($addr in @("B4","C4","D4","L4","M4","N4","O4")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("B4").Value = "23.00"
$ws.Range("C4").Value = "-220"
$ws.Range("D4").Value = "-125"
$ws.Range("E4").Value = 29.6
$ws.Range("F4").Value = 25.3
$ws.Range("G4").Value = "24.30+"
$ws.Range("H4").Value = "25.82+"
$ws.Range("I4").Value = "23.13+"
$ws.Range("J4").Value = "- (887668.60)"
$ws.Range("K4").Value = "-"
$ws.Range("L4").Value = "-0.58"
$ws.Range("M4").Value = "25.00"
$ws.Range("N4").Value = "0.63"
$ws.Range("O4").Value = "20.71"
$ws.Range("P4").Value = "+--+++---+++"
$ws.Range("Q4").Value = 7
$ws.Range("R4").Value = 5

# 6) New row 5: 1264.TWO, full indicator row (OTC listing).
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "1264.TWO"
foreach ($addr in @("B5","C5","D5","L5","M5","N5","O5")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("B5").Value = "295.00"
$ws.Range("C5").Value = "-1"
$ws.Range("D5").Value = "0"
$ws.Range("E5").Value = 294.5
$ws.Range("F5").Value = 289
$ws.Range("G5").Value = "294.77-"
$ws.Range("H5").Value = "292.52-"
$ws.Range("I5").Value = "285.47-"
$ws.Range("J5").Value = "- (820.05)"
$ws.Range("K5").Value = "+"
$ws.Range("L5").Value = "0.66"
$ws.Range("M5").Value = "37.50"
$ws.Range("N5").Value = "2.59"
$ws.Range("O5").Value = "76.00"
$ws.Range("P5").Value = "+------+++++"
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 6
